{"js": "// Ordered list of replacement texts: index 0 is the date title paragraph,\n// indices 1..100 are the 100 arithmetic-expression table-cell paragraphs,\n// in document (reading) order - this matches the order\n// `context.document.body.paragraphs` enumerates them in.\nconst newValues = [\n  \"2025-02-07 Friday\",\n  \"31+52=\",\n  \"13+19=\",\n  \"15+39=\",\n  \"92+3=\",\n  \"28+24=\",\n  \"87+2=\",\n  \"56-50=\",\n  \"24+43=\",\n  \"97-61=\",\n  \"76-48=\",\n  \"26+23=\",\n  \"47-43=\",\n  \"97-84=\",\n  \"10+63=\",\n  \"46-16=\",\n  \"73+19=\",\n  \"50-27=\",\n  \"9+58=\",\n  \"41+10=\",\n  \"79+17=\",\n  \"14+7=\",\n  \"61+32=\",\n  \"44+38=\",\n  \"25+38=\",\n  \"99-9=\",\n  \"27+48=\",\n  \"68-6=\",\n  \"47-26=\",\n  \"73-67=\",\n  \"74-56=\",\n  \"82-80=\",\n  \"26+2=\",\n  \"12+10=\",\n  \"86-23=\",\n  \"41-19=\",\n  \"55-10=\",\n  \"75+23=\",\n  \"33+46=\",\n  \"9+40=\",\n  \"65+12=\",\n  \"25-18=\",\n  \"10+34=\",\n  \"98-18=\",\n  \"89-85=\",\n  \"29+69=\",\n  \"9+36=\",\n  \"72-27=\",\n  \"88-70=\",\n  \"6-3=\",\n  \"6+52=\",\n  \"58+7=\",\n  \"72+9=\",\n  \"32+11=\",\n  \"64+30=\",\n  \"27-13=\",\n  \"91-16=\",\n  \"11+18=\",\n  \"78-60=\",\n  \"10+24=\",\n  \"12+41=\",\n  \"28-7=\",\n  \"26+51=\",\n  \"45+46=\",\n  \"59+12=\",\n  \"9+21=\",\n  \"93+1=\",\n  \"19+7=\",\n  \"53-35=\",\n  \"47+33=\",\n  \"86-57=\",\n  \"50+21=\",\n  \"97-44=\",\n  \"42-29=\",\n  \"34+24=\",\n  \"52-2=\",\n  \"39+28=\",\n  \"74-10=\",\n  \"21+57=\",\n  \"24+16=\",\n  \"27+17=\",\n  \"82+1=\",\n  \"63-57=\",\n  \"40-25=\",\n  \"22+67=\",\n  \"12+14=\",\n  \"0+90=\",\n  \"91-31=\",\n  \"16+53=\",\n  \"45-18=\",\n  \"97-52=\",\n  \"9+43=\",\n  \"93-29=\",\n  \"21+6=\",\n  \"22-7=\",\n  \"58-24=\",\n  \"89-70=\",\n  \"52+40=\",\n  \"26-12=\",\n  \"68-3=\",\n  \"83-5=\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(`Expected ${newValues.length} paragraphs, found ${paragraphs.items.length}`);\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date title (first paragraph of the document)\n$d = $word.ActiveDocument\n$d.Paragraphs.Item(1).Range.Text = '2025-02-07 Friday'\n\n# Update each of the 100 arithmetic-expression table cells, in row-major\n# (reading) order, matching the order the cells appear in the document.\n$newValues = @(\n    '31+52=',\n    '13+19=',\n    '15+39=',\n    '92+3=',\n    '28+24=',\n    '87+2=',\n    '56-50=',\n    '24+43=',\n    '97-61=',\n    '76-48=',\n    '26+23=',\n    '47-43=',\n    '97-84=',\n    '10+63=',\n    '46-16=',\n    '73+19=',\n    '50-27=',\n    '9+58=',\n    '41+10=',\n    '79+17=',\n    '14+7=',\n    '61+32=',\n    '44+38=',\n    '25+38=',\n    '99-9=',\n    '27+48=',\n    '68-6=',\n    '47-26=',\n    '73-67=',\n    '74-56=',\n    '82-80=',\n    '26+2=',\n    '12+10=',\n    '86-23=',\n    '41-19=',\n    '55-10=',\n    '75+23=',\n    '33+46=',\n    '9+40=',\n    '65+12=',\n    '25-18=',\n    '10+34=',\n    '98-18=',\n    '89-85=',\n    '29+69=',\n    '9+36=',\n    '72-27=',\n    '88-70=',\n    '6-3=',\n    '6+52=',\n    '58+7=',\n    '72+9=',\n    '32+11=',\n    '64+30=',\n    '27-13=',\n    '91-16=',\n    '11+18=',\n    '78-60=',\n    '10+24=',\n    '12+41=',\n    '28-7=',\n    '26+51=',\n    '45+46=',\n    '59+12=',\n    '9+21=',\n    '93+1=',\n    '19+7=',\n    '53-35=',\n    '47+33=',\n    '86-57=',\n    '50+21=',\n    '97-44=',\n    '42-29=',\n    '34+24=',\n    '52-2=',\n    '39+28=',\n    '74-10=',\n    '21+57=',\n    '24+16=',\n    '27+17=',\n    '82+1=',\n    '63-57=',\n    '40-25=',\n    '22+67=',\n    '12+14=',\n    '0+90=',\n    '91-31=',\n    '16+53=',\n    '45-18=',\n    '97-52=',\n    '9+43=',\n    '93-29=',\n    '21+6=',\n    '22-7=',\n    '58-24=',\n    '89-70=',\n    '52+40=',\n    '26-12=',\n    '68-3=',\n    '83-5='\n)\n\n$t = $d.Tables.Item(1)\nif ($t.Range.Cells.Count -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) cells, found $($t.Range.Cells.Count)\"\n}\n\n$i = 0\nforeach ($cell in $t.Range.Cells) {\n    $cell.Range.Text = $newValues[$i]\n    $i = $i + 1\n}\n"}
